$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to stay Text so the stored value matches the original string cells.
$textCells = @('D5', 'D8', 'D11', 'D14', 'D18', 'D20', 'D21', 'D22', 'D25', 'D27', 'D28', 'D30', 'D33', 'D36', 'D40', 'D42', 'D43', 'D48')
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range('D2').Value = '35.192.12'
$ws.Range('E2').Value = '  +1.24%  '
$ws.Range('D3').Value = '1.857.77'
$ws.Range('E3').Value = '  +1.67%  '
$ws.Range('E4').Value = '  +0.48%  '
$ws.Range('D5').Value = '239.13'
$ws.Range('E5').Value = '  +3.81%  '
$ws.Range('E6').Value = '  +0.78%  '
$ws.Range('E7').Value = '  +0.48%  '
$ws.Range('D8').Value = '42.07'
$ws.Range('E8').Value = '  +6.62%  '
$ws.Range('E9').Value = '  +1.11%  '
$ws.Range('D11').Value = '0.0988'
$ws.Range('E11').Value = '  +0.21%  '
$ws.Range('D12').Value = '2.126.84'
$ws.Range('E12').Value = '  +1.68%  '
$ws.Range('E13').Value = '  +1.98%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = '0.676'
$ws.Range('E14').Value = '  +1.11%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '1.829.01'
$ws.Range('E15').Value = '  +0.09%  '
$ws.Range('E16').Value = '  +2.51%  '
$ws.Range('D17').Value = '35.175.18'
$ws.Range('E17').Value = '  +1.21%  '
$ws.Range('D18').Value = '69.86'
$ws.Range('E18').Value = '  +0.57%  '
$ws.Range('E19').Value = '  +1.43%  '
$ws.Range('D20').Value = '240.75'
$ws.Range('E20').Value = '  +0.21%  '
$ws.Range('D21').Value = '12.21'
$ws.Range('E21').Value = '  +0.95%  '
$ws.Range('D22').Value = '4.75'
$ws.Range('E22').Value = '  +1.85%  '
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('E24').Value = '  +0.95%  '
$ws.Range('D25').Value = '169.57'
$ws.Range('E25').Value = '  -1.10%  '
$ws.Range('E26').Value = '  +26.26%  '
$ws.Range('D27').Value = '8.01'
$ws.Range('E27').Value = '  +3.60%  '
$ws.Range('D28').Value = '17.66'
$ws.Range('E28').Value = '  +1.99%  '
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('D30').Value = '0.0562'
$ws.Range('E30').Value = '  +2.35%  '
$ws.Range('E31').Value = '  +0.45%  '
$ws.Range('E32').Value = '  +2.15%  '
$ws.Range('D33').Value = '1.81'
$ws.Range('E33').Value = '  +26.21%  '
$ws.Range('E34').Value = '  +2.15%  '
$ws.Range('E35').Value = '  +11.07%  '
$ws.Range('D36').Value = '0.817'
$ws.Range('E36').Value = '  +17.34%  '
$ws.Range('E37').Value = '  +7.55%  '
$ws.Range('E38').Value = '  +4.46%  '
$ws.Range('E39').Value = '  +4.50%  '
$ws.Range('D40').Value = '89.89'
$ws.Range('E40').Value = '  -1.26%  '
$ws.Range('D41').Value = '1.346.44'
$ws.Range('E41').Value = '  +0.68%  '
$ws.Range('D42').Value = '0.0601'
$ws.Range('E42').Value = '  +15.09%  '
$ws.Range('D43').Value = '14.98'
$ws.Range('E43').Value = '  +3.42%  '
$ws.Range('E44').Value = '  +2.86%  '
$ws.Range('E45').Value = '  +0.26%  '
$ws.Range('E46').Value = '  +43.84%  '
$ws.Range('E47').Value = '  -0.43%  '
$ws.Range('D48').Value = '6.58'
$ws.Range('E48').Value = '  +5.39%  '
$ws.Range('D49').Value = '2.045.41'
$ws.Range('E49').Value = '  +1.92%  '
$ws.Range('E50').Value = '  +1.46%  '
$ws.Range('E51').Value = '  +0.50%  '

# Drop the temporary Text number format again so no stray style survives the edit.
foreach ($addr in $textCells) { $ws.Range($addr).ClearFormats() }
